$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 values changed by the edit. A4/B4/D4 look numeric, so without any
# precaution Excel would silently coerce them to doubles (and the 19-digit
# ID in A4 would lose precision / turn into scientific notation). Prefixing
# with a single quote forces Excel to store them as text, matching the
# original inlineStr cells. ClearFormats() afterwards removes the implicit
# "Text" number-format/style that typing a quote-prefixed value applies, so
# the cell keeps its original (unstyled) look.
$ws.Range("A4").Value = "'7264046483537334765"
$ws.Range("A4").ClearFormats()

$ws.Range("B4").Value = "'29"
$ws.Range("B4").ClearFormats()

$ws.Range("D4").Value = "'2022"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").Value = "Ein Tag vor Berkos Bday"
